# regen sval data to filter save games
#
# The underlying s_val dataset was regenerated (filtering out "save game"
# rows before computing the per-stat averages), which changes the TB/d2S/K/IP
# values (columns B:E) for every row, and consequently the rolling "sum"
# column (G = TB + d2S + K + IP). The "Win" column (F) and the date labels
# (column A) are unaffected.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New TB, d2S, K, IP, sum values per row (row number => B,C,D,E,G)
$rows = @{
    2  = @(3.182878228561681,   1.65323645889881,    3.082599426703578,   6.48142807727062,  14.40014219143469)
    3  = @(0.3464964993005633,  0.3375848360084654,  0.1529057820181812,  0.4998867070740569, 1.336873824401267)
    4  = @(3.182878228561681,   1.65323645889881,    3.082599426703578,   0.4998867070740569, 8.418600821238126)
    5  = @(0.7287194209349384,  1.65323645889881,    0.7127328510149897,  0.4998867070740569, 3.594575437922795)
    6  = @(0.1554434735375247,  0.3375848360084654,  3.082599426703578,   6.48142807727062,  10.05705581352019)
    7  = @(3.182878228561681,   1.65323645889881,    0.7127328510149897,  6.48142807727062,  12.0302756157461)
    8  = @(0.02258322285507441, 1.65323645889881,    0.1529057820181812,  0.4998867070740569, 2.328612170846122)
    9  = @(3.182878228561681,   9.226618575922256,   16.98373111632243,   6.48142807727062,  35.87465599807698)
    10 = @(3.182878228561681,   87981.0709163148,    3.082599426703578,   6.48142807727062,  87993.81782204733)
    11 = @(1.505614041169197,   1.65323645889881,    3.082599426703578,   0.4998867070740569, 6.741336633845642)
    12 = @(0.1554434735375247,  0.05231270169004087, 0.7127328510149897,  0.4998867070740569, 1.420375733316612)
    13 = @(3.182878228561681,   1.65323645889881,    157.8057217802531,   6.48142807727062, 169.1232645449842)
    14 = @(1.505614041169197,   1.65323645889881,    0.7127328510149897,  0.4998867070740569, 4.371470058157054)
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Range("B$r").Value = $vals[0]
    $ws.Range("C$r").Value = $vals[1]
    $ws.Range("D$r").Value = $vals[2]
    $ws.Range("E$r").Value = $vals[3]
    $ws.Range("G$r").Value = $vals[4]
}
